$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 'D2' '40.654.11'
Set-TextValue 'E2' '  -2.33%  '
Set-TextValue 'D3' '2.372.31'
Set-TextValue 'E3' '  -4.09%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '310.97'
Set-TextValue 'E5' '  -2.50%  '
Set-TextValue 'D6' '86.69'
Set-TextValue 'E6' '  -6.10%  '
Set-TextValue 'D7' '0.530'
Set-TextValue 'E7' '  -3.96%  '
Set-TextValue 'E9' '  -4.50%  '
Set-TextValue 'D10' '0.0840'
Set-TextValue 'E10' '  -3.11%  '
Set-TextValue 'D11' '30.39'
Set-TextValue 'E11' '  -8.28%  '
Set-TextValue 'D12' '0.110'
Set-TextValue 'E12' '  -0.58%  '
Set-TextValue 'D13' '2.737.48'
Set-TextValue 'E13' '  -4.12%  '
Set-TextValue 'D14' '6.52'
Set-TextValue 'E14' '  -5.51%  '
Set-TextValue 'D15' '14.96'
Set-TextValue 'E15' '  -3.41%  '
Set-TextValue 'D16' '2.389.95'
Set-TextValue 'E16' '  -2.97%  '
Set-TextValue 'D17' '0.757'
Set-TextValue 'E17' '  -4.22%  '
Set-TextValue 'D18' '40.588.47'
Set-TextValue 'E18' '  -2.41%  '
Set-TextValue 'D19' '0.0₃0910'
Set-TextValue 'E19' '  -3.67%  '
Set-TextValue 'E20' '  -4.99%  '
Set-TextValue 'D21' '68.47'
Set-TextValue 'E21' '  -3.14%  '
Set-TextValue 'D22' '10.75'
Set-TextValue 'E22' '  -4.67%  '
Set-TextValue 'D23' '235.24'
Set-TextValue 'E23' '  -2.24%  '
Set-TextValue 'D24' '2.57'
Set-TextValue 'E24' '  -6.27%  '
Set-TextValue 'E25' '  +0.09%  '
Set-TextValue 'D26' '1.81'
Set-TextValue 'E26' '  -7.48%  '
Set-TextValue 'D27' '23.81'
Set-TextValue 'E27' '  -3.99%  '
Set-TextValue 'E28' '  -0.47%  '
Set-TextValue 'E29' '  -4.74%  '
Set-TextValue 'D30' '34.23'
Set-TextValue 'E30' '  -6.95%  '
Set-TextValue 'D31' '153.86'
Set-TextValue 'E31' '  -2.06%  '
Set-TextValue 'E32' '  -0.05%  '
Set-TextValue 'D33' '5.18'
Set-TextValue 'E33' '  -5.26%  '
Set-TextValue 'D34' '0.0728'
Set-TextValue 'E34' '  -4.81%  '
Set-TextValue 'E35' '  -5.90%  '
Set-TextValue 'D36' '0.114'
Set-TextValue 'E36' '  -2.34%  '
Set-TextValue 'D37' '15.98'
Set-TextValue 'E37' '  -7.72%  '
Set-TextValue 'E38' '  -4.11%  '
Set-TextValue 'D39' '2.76'
Set-TextValue 'E39' '  -4.65%  '
Set-TextValue 'D40' '1.70'
Set-TextValue 'E40' '  -8.10%  '
Set-TextValue 'E41' '  -3.38%  '
Set-TextValue 'D42' '2.41'
Set-TextValue 'E42' '  -3.07%  '
Set-TextValue 'D43' '1.957.27'
Set-TextValue 'E43' '  -1.46%  '
Set-TextValue 'D44' '0.0269'
Set-TextValue 'E44' '  -4.97%  '
Set-TextValue 'D45' '17.62'
Set-TextValue 'E45' '  -6.71%  '
Set-TextValue 'D46' '9.40'
Set-TextValue 'E46' '  -0.70%  '
Set-TextValue 'E47' '  -9.74%  '
Set-TextValue 'D48' '2.592.80'
Set-TextValue 'E48' '  -4.39%  '
Set-TextValue 'D49' '93.01'
Set-TextValue 'E49' '  -5.07%  '
Set-TextValue 'D50' '71.78'
Set-TextValue 'E50' '  -5.20%  '
Set-TextValue 'D51' '50.16'
Set-TextValue 'E51' '  -4.30%  '
